$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) column C for every data row (2 through 410)
# is being bumped from serial date 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C410").Value = 45175
